$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 09:15"

# Update row 8 values (Benzina Albert Modrice)
$ws.Range("B8").Value = 37.5
$ws.Range("C8").Value = 36.5

# D8 switches from a numeric delta to a text label; force text so Excel
# doesn't reinterpret "+1.0" as the number 1, then drop the resulting
# number-format style so the cell stays unstyled like the rest of the row.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "+1.0"
$ws.Range("D8").ClearFormats()

# E8 switches from a numeric date serial (formatted cell) to a plain text
# timestamp string; clear the old date number format so the cell reverts
# to default (unstyled) and holds the literal text instead of a serial.
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = "2022-02-16 09:15:54"
